$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.307.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4757"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06449"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.848.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07414"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.009"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6347"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.279.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  -3.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007392"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.099.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.121"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.038"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.291"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("E28").Value = "  -4.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1029"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.57%  "
$ws.Range("E30").Value = "  -5.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.246"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.919"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.151"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7290"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9987"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01962"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.41%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9077"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.990"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4126"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.573"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.76%  "
$ws.Range("E46").Value = "  -4.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  -4.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.852"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.403"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.68%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05610"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
